$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: pup F juvenile(?) values updated; J3 formula divisor changed ---
$ws.Range("F3").Value = 0.14699999999999999
$ws.Range("J3").Formula = "=(I3/923.4)*1000"

# --- Row 7 ---
$ws.Range("F7").Value = 0.49099999999999999
$ws.Range("J7").Formula = "=(I7/223.2)*1000"

# --- Row 9 ---
$ws.Range("E9").Value = 0.153

# --- Row 13 ---
$ws.Range("E13").Value = 0.093

# --- Row 14: new unit-conversion label cell ---
$ws.Range("I14").Value = "mlO2/kg*min"

# --- Row 15 ---
$ws.Range("I15").Value = 1.804
$ws.Range("J15").Formula = "=I15*20.08/1000"

# --- Row 17: F17 becomes a literal value, I17 updated, J17 new formula ---
$ws.Range("F17").Value = 0.126
$ws.Range("I17").Value = 6.2619999999999996
$ws.Range("J17").Formula = "=I17*20.08/1000"

# --- Row 19 ---
$ws.Range("I19").Value = 5.3769999999999998
$ws.Range("J19").Formula = "=I19*20.08/1000"

# --- Row 20 (blank separator row gains a J formula) ---
$ws.Range("J20").Formula = "=I20*20.08/1000"

# --- Row 21 ---
$ws.Range("J21").Formula = "=I21*20.08/1000"

# --- Row 27 ---
$ws.Range("I27").Value = 0.90468231877835403
$ws.Range("J27").Formula = "=I27*20.08/1000"

# --- Row 29 ---
$ws.Range("I29").Value = 2.9506455721203002
$ws.Range("J29").Formula = "=I29*20.08/1000"

# --- Row 31 ---
$ws.Range("I31").Value = 1.6833771689994701
$ws.Range("J31").Formula = "=I31*20.08/1000"

# --- Row 33 ---
$ws.Range("E33").Value = 0.085

# --- Row 37 ---
$ws.Range("E37").Value = 0.091
$ws.Range("F37").Value = 0.107

# --- Row 39 ---
$ws.Range("E39").Value = 0.153

# --- Row 43 ---
$ws.Range("E43").Value = 0.093

# --- Row 45 ---
$ws.Range("F45").Value = 0.14699999999999999

# --- Row 47 ---
$ws.Range("F47").Value = 0.85399999999999998

# --- Row 49 ---
$ws.Range("F49").Value = 0.49099999999999999

# --- Selection moved to E43 (last active cell in the sheet) ---
$ws.Range("E43").Select()
